# cages.xlsx edit:
#   - fix row 2 (cage #1) values, and change its material to Metal
#   - append two more cages (rows 10 and 11)
#   - re-sort the cage list by CageNumber (column A), including the new rows
#   - leave the active selection on F9 (next empty cell after the table)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- correct cage #1 (row 2) ---
$ws.Range("B2").Value = 19
$ws.Range("C2").Value = 20
$ws.Range("D2").Value = 18
$ws.Range("E2").Value = "Metal"

# --- add cage #9 (row 10) ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 30
$ws.Range("D10").Value = 40
$ws.Range("E10").Value = "Wood"

# --- add cage #10 (row 11) ---
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 15
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 25
$ws.Range("E11").Value = "Plastic"

# --- sort the whole table (with header) by CageNumber ascending ---
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2"))
$sortObj.SetRange($ws.Range("A1:E11"))
$sortObj.Header = 1
$sortObj.Apply()

# --- move the selection to the first empty cell below/right of the table ---
$null = $ws.Range("F9").Select()
